$d = $word.ActiveDocument

# 1. "possible" -> "possibly"
$d.Content.Find.Execute("could possible ruin", $true, $false, $false, $false, $false, $true, 1, $false, "could possibly ruin", 2) | Out-Null

# 2. Remove the "math-oriented / pre-image attack" sentence from the MD5 explanation paragraph
$d.Content.Find.Execute(" For our more math-oriented people this is known as a " + [char]8220 + "pre-image" + [char]8221 + " attack, since the input of a function/mapping is called the pre-image. Now", $true, $false, $true, $false, $false, $true, 1, $false, " Now", 2) | Out-Null

# 3. Remove the "pre-image" qualifier before "values"
$d.Content.Find.Execute("associated pre-image values", $true, $false, $false, $false, $false, $true, 1, $false, "associated values", 2) | Out-Null

# 4. Remove the trailing "hashcat/crunch" paragraph text, keeping only the image in that paragraph
$d.Content.Find.Execute("If you feel adventurous*careful.", $true, $false, $true, $false, $false, $true, 1, $false, "", 2) | Out-Null
